$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) values change from serial date 45221 (2023-10-22)
# to serial date 45224 (2023-10-25) for all data rows (2 through 20).
for ($row = 2; $row -le 20; $row++) {
    $ws.Cells.Item($row, 3).Value = 45224
}
